$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update VisitsPerMonth (column B) values
$ws.Range("B2").Value = 15
$ws.Range("B3").Value = 15
$ws.Range("B4").Value = 15
$ws.Range("B5").Value = 15
$ws.Range("B6").Value = 15
$ws.Range("B9").Value = 29
$ws.Range("B10").Value = 30
$ws.Range("B11").Value = 15
$ws.Range("B12").Value = 15
$ws.Range("B13").Value = 15
$ws.Range("B14").Value = 15

# Update ExpectedRank (column D) values - "Gold" becomes "Standard"
$ws.Range("D4").Value = "Standard"
$ws.Range("D5").Value = "Standard"
$ws.Range("D6").Value = "Standard"
$ws.Range("D13").Value = "Standard"
$ws.Range("D14").Value = "Standard"

# Update selection to reflect the reviewed/edited range
$ws.Range("E2:F14").Select()
